$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: location -> y, lon/lat/altitude updated
$ws.Range("B2").Value = "y"
$ws.Range("C2").Value = 112.15
$ws.Range("D2").Value = 22
$ws.Range("E2").Value = 50

# Update row 3: location -> z, lon/lat/altitude updated
$ws.Range("B3").Value = "z"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -30
$ws.Range("E3").Value = 15

# Remove row 4 entirely (was the "3 | z | 116 | 40 | 1000" row)
$ws.Range("A4:E4").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Update the selection to match the saved view state
$ws.Range("B5").Select()
